# Update NPC1 ExAC data on the "Ranked_Combined" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ranked_Combined")

# Make sure this sheet stays the active / selected tab, as it was before.
$ws.Activate() | Out-Null

# Try to restore the sheet-tab split ratio (cosmetic; may not round-trip).
$win = $excel.ActiveWindow
$win.TabRatio = 0.992

function Set-Lit($cellRef, $term) {
    $url = "http://www.ncbi.nlm.nih.gov/sites/entrez?db=pubmed&cmd=search&term=" + $term
    $ws.Range($cellRef).Value = $url
}

# Row 4 (V889M)
$ws.Range("C4").Value = 1
Set-Lit "D4" "11182931"

# Row 5 (P1007A)
$ws.Range("C5").Value = -1
Set-Lit "D5" "11754101"

# Row 6 (R1186H)
$ws.Range("C6").Value = 0
Set-Lit "D6" "11349231"
Set-Lit "E6" "11333381"
Set-Lit "G6" "16126423"

# Row 8 (T1205K)
$ws.Range("C8").Value = -1
Set-Lit "D8" "12955717"

# Row 9 (V378A)
$ws.Range("C9").Value = 0
Set-Lit "D9" "11333381"

# Row 10 (R404Q)
$ws.Range("C10").Value = 1
Set-Lit "D10" "11349231"
Set-Lit "E10" "11333381"
Set-Lit "F10" "11545687"

# Row 11 (I1061T)
$ws.Range("C11").Value = -1
Set-Lit "D11" "11754101"
Set-Lit "E11" "10521290"
Set-Lit "F11" "10521297"
Set-Lit "G11" "11182931"
Set-Lit "H11" "11349231"
Set-Lit "I11" "11333381"
Set-Lit "J11" "11479732"
Set-Lit "K11" "12401890"
Set-Lit "L11" "16098014"
Set-Lit "M11" "16126423"

# Row 21 (P434S)
$ws.Range("C21").Value = -1
Set-Lit "D21" "12955717"

# Restore the previously-selected cell on this sheet.
$ws.Range("D11").Select() | Out-Null
